# Sieving Policy document update
#
# draft cluster_cavgs_selection commander to support updated sieving
# policy. Parallelization of affinity propagation clustering.
#
# Mutates $word.ActiveDocument from the "before" state to the "after"
# state described by the commit's OOXML diff.

$d = $word.ActiveDocument

function Insert-TextAtParaEnd($doc, $para, $text) {
    # Inserts $text just before the paragraph mark of $para, returning a
    # Range covering exactly the inserted text (so callers can format it).
    $rr = $para.Range
    $startPos = $rr.End - 1
    $ins = $doc.Range($startPos, $startPos)
    $ins.InsertAfter($text)
    $endPos = $startPos + $text.Length
    return $doc.Range($startPos, $endPos)
}

function Clear-ParaText($doc, $para) {
    $rr = $para.Range
    $full = $doc.Range($rr.Start, $rr.End - 1)
    $full.Text = ""
}

# ---------------------------------------------------------------------
# 1) Remove the blank paragraph between the title and the first bullet.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# ---------------------------------------------------------------------
# 2) Bullet 1 "Chunks of 2D are done in 100 micrograph lumps with 100
#    classes per chunk" -> symbolic N (subscript M) / N (subscript C)
#    counts instead of the literal "100" values.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(2)
Clear-ParaText $d $p1
Insert-TextAtParaEnd $d $p1 "Chunks of 2D are done in " | Out-Null
Insert-TextAtParaEnd $d $p1 "N" | Out-Null
$sub1 = Insert-TextAtParaEnd $d $p1 "M"
$sub1.Font.Subscript = -1
Insert-TextAtParaEnd $d $p1 " micrograph lumps with " | Out-Null
Insert-TextAtParaEnd $d $p1 "N" | Out-Null
$sub2 = Insert-TextAtParaEnd $d $p1 "C"
$sub2.Font.Subscript = -1
Insert-TextAtParaEnd $d $p1 " classes per chunk" | Out-Null

# ---------------------------------------------------------------------
# 3) Bullet 2 "Clustering of 2D classes is done..." is unchanged.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 4) Bullet 3 "Sets are combined up to 10 sets..." -> reworded/expanded
#    explanation of the merging/clustering behaviour.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(4)
Clear-ParaText $d $p3
Insert-TextAtParaEnd $d $p3 "Sets are combined" | Out-Null
Insert-TextAtParaEnd $d $p3 ". Until the set of combined class averages reaches 1000 members, the clustering of the merged set is done without reference to any previous clustering solution, thus " | Out-Null
Insert-TextAtParaEnd $d $p3 "replacing the prior solution with each integration of a new " | Out-Null
Insert-TextAtParaEnd $d $p3 "chunk" | Out-Null
Insert-TextAtParaEnd $d $p3 " to create a " | Out-Null
Insert-TextAtParaEnd $d $p3 "combined " | Out-Null
Insert-TextAtParaEnd $d $p3 "set with " | Out-Null
Insert-TextAtParaEnd $d $p3 "a minimum of " | Out-Null
Insert-TextAtParaEnd $d $p3 "1000 classes in it" | Out-Null
Insert-TextAtParaEnd $d $p3 "." | Out-Null

# ---------------------------------------------------------------------
# 5) Bullet 4 "At this point the user is prompted..." is unchanged.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 6) Bullet 5 "The good classes and bad classes are then re-clustered
#    separately..." -> drop the bracketed "[represent_selection]" and
#    instead name the new commander "cluster_cavgs_selection" in a full
#    sentence, keeping the spell-check wrapper around the identifier.
# ---------------------------------------------------------------------

# Drop the trailing " [" that introduced the old bracketed reference.
$d.Content.Find.Execute("matching [", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "matching", 2) | Out-Null

# Insert the new explanatory sentence right after "...further matching",
# before the (still proofErr-wrapped) identifier run.
$fr = $d.Content
$fr.Find.Execute("further matching", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$newSentencePos = $fr.End
$newSentenceIns = $d.Range($newSentencePos, $newSentencePos)
$newSentenceIns.InsertAfter(". The program that implements this " + [char]0x201C + `
                             "medoid-for-matching-generation" + [char]0x201D + " is called ")

# Rename the identifier itself (stays wrapped by the same proofErr tags).
$d.Content.Find.Execute("represent_selection", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "cluster_cavgs_selection", 2) | Out-Null

# Turn the trailing "]" run into a closing "." (leave the identifier run
# and its spellEnd proofErr marker untouched).
$fb = $d.Content
$fb.Find.Execute("cluster_cavgs_selection]", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$bracketRange = $d.Range($fb.End - 1, $fb.End)
$bracketRange.Text = "."

# ---------------------------------------------------------------------
# 7) Bullet 6 "We then proceed by simply matching incoming chunks..."
#    gains a trailing full stop after the closing parenthesis.
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(7)
Insert-TextAtParaEnd $d $p6 "." | Out-Null
